# Apply "dSF" (column F) corrections for specific rows, per commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    6  = -1
    8  = 0
    9  = -2
    10 = -1
    17 = 2
    19 = 1
    23 = 3
    30 = 2
    32 = 1
    39 = -1
    42 = 0
    51 = 3
    53 = 1
    54 = 7
    57 = 1
    59 = -6
    60 = 2
    61 = -5
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
